$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 152; existing rows 152-160 shift down to 153-161.
$ws.Rows(152).Insert()

# Populate the newly inserted row 152 with the new weekly price entry.
$ws.Cells.Item(152, 1).Value = 7
$ws.Cells.Item(152, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(152, 3).Value = "Ñuble"
$ws.Cells.Item(152, 4).Value = 44753
$ws.Cells.Item(152, 5).Value = 16
$ws.Cells.Item(152, 6).Value = 100112045
$ws.Cells.Item(152, 7).Value = "Zapallo"
$ws.Cells.Item(152, 8).Value = "Camote"
$ws.Cells.Item(152, 9).Value = "1a (guarda)"
$ws.Cells.Item(152, 10).Value = 200
$ws.Cells.Item(152, 11).Value = 500
$ws.Cells.Item(152, 12).Value = 550
$ws.Cells.Item(152, 13).Value = 525
$ws.Cells.Item(152, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(152, 15).Value = "Región del Maule"
$ws.Cells.Item(152, 16).Value = 525
$ws.Cells.Item(152, 17).Value = 1
$ws.Cells.Item(152, 18).Value = "Hortaliza"
